# Update the build timestamp embedded in the version string throughout the
# workbook, from "February 03 2026 17.29.55 EST" to "February 03 2026 18.05.36 EST".

$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

# --- "About" sheet ---
$wsAbout = $wb.Worksheets.Item("About")

$rA2 = $wsAbout.Range("A2")
$a2 = $rA2.Value()
$rA2.Value = $a2.Replace($oldStamp, $newStamp)

$rA6 = $wsAbout.Range("A6")
$a6 = $rA6.Value()
$rA6.Value = $a6.Replace($oldStamp, $newStamp)

# --- "Boundaries and methane sources" sheet ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 11; $row++) {
    $cell = $wsData.Cells.Item($row, 19)  # column S = build_version
    $val = $cell.Value()
    $cell.Value = $val.Replace($oldStamp, $newStamp)
}
